$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while forcing text storage so that
# numeric-looking strings (e.g. "416.11", "1.00") are preserved verbatim
# instead of being auto-converted into floating point numbers by Excel.
# We briefly apply a text NumberFormat, assign the value, then restore the
# cell style back to "Normal" so no residual style index is left on the cell
# (matching the original workbook, where these cells carry no explicit style).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.240.43"
$ws.Range("E2").Value = "  +7.11%  "
Set-TextValue $ws.Range("D3") "3.591.90"
$ws.Range("E3").Value = "  +3.60%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue $ws.Range("D5") "416.11"
$ws.Range("E5").Value = "  +0.58%  "
Set-TextValue $ws.Range("D6") "129.82"
$ws.Range("E6").Value = "  -0.26%  "
Set-TextValue $ws.Range("D7") "0.654"
$ws.Range("E7").Value = "  +4.17%  "
Set-TextValue $ws.Range("D8") "3.582.57"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("E9").Value = "  -0.05%  "
Set-TextValue $ws.Range("D10") "0.779"
$ws.Range("E10").Value = "  +7.22%  "
$ws.Range("E11").Value = "  +17.43%  "
Set-TextValue $ws.Range("D12") "0.0000339"
$ws.Range("E12").Value = "  +54.47%  "
Set-TextValue $ws.Range("D13") "42.62"
$ws.Range("E13").Value = "  +0.19%  "
Set-TextValue $ws.Range("D14") "9.92"
$ws.Range("E14").Value = "  +3.30%  "
Set-TextValue $ws.Range("D15") "4.168.63"
$ws.Range("E15").Value = "  +3.76%  "
$ws.Range("E16").Value = "  -0.23%  "
Set-TextValue $ws.Range("D17") "20.33"
$ws.Range("E17").Value = "  -1.14%  "
Set-TextValue $ws.Range("D18") "3.609.27"
$ws.Range("E18").Value = "  +4.40%  "
$ws.Range("E19").Value = "  +5.39%  "
Set-TextValue $ws.Range("D20") "67.106.33"
$ws.Range("E20").Value = "  +7.05%  "
$ws.Range("E21").Value = "  -3.17%  "
Set-TextValue $ws.Range("D22") "452.96"
$ws.Range("E22").Value = "  -1.99%  "
Set-TextValue $ws.Range("D23") "89.68"
$ws.Range("E23").Value = "  -0.98%  "
Set-TextValue $ws.Range("D24") "3.17"
$ws.Range("E24").Value = "  -3.27%  "
Set-TextValue $ws.Range("D25") "13.16"
$ws.Range("E25").Value = "  -1.11%  "
Set-TextValue $ws.Range("D26") "3.36"
$ws.Range("E26").Value = "  +1.57%  "
Set-TextValue $ws.Range("D27") "10.06"
$ws.Range("E27").Value = "  -6.26%  "
Set-TextValue $ws.Range("D28") "35.13"
$ws.Range("E28").Value = "  +5.15%  "
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("E31").Value = "  +3.95%  "
Set-TextValue $ws.Range("D32") "0.118"
$ws.Range("E32").Value = "  +5.07%  "
Set-TextValue $ws.Range("D33") "7.40"
$ws.Range("E33").Value = "  -2.39%  "
Set-TextValue $ws.Range("D34") "0.162"
$ws.Range("E34").Value = "  -3.21%  "
Set-TextValue $ws.Range("D35") "40.70"
$ws.Range("E35").Value = "  -0.16%  "
Set-TextValue $ws.Range("D36") "1.00"
$ws.Range("E36").Value = "  +0.05%  "
Set-TextValue $ws.Range("D37") "56.72"
$ws.Range("E37").Value = "  -3.04%  "
Set-TextValue $ws.Range("D38") "0.0495"
$ws.Range("E38").Value = "  +0.90%  "
Set-TextValue $ws.Range("D39") "0.0₃0745"
$ws.Range("E39").Value = "  +33.83%  "
Set-TextValue $ws.Range("D40") "0.148"
$ws.Range("E40").Value = "  +10.30%  "
Set-TextValue $ws.Range("D41") "0.999"
$ws.Range("E41").Value = "  +0.02%  "
Set-TextValue $ws.Range("D42") "3.06"
$ws.Range("E42").Value = "  -1.02%  "
Set-TextValue $ws.Range("D43") "149.92"
$ws.Range("E43").Value = "  +1.67%  "
Set-TextValue $ws.Range("D44") "2.75"
$ws.Range("E44").Value = "  +2.39%  "
# Row 45: coin identity/rank changed (values swapped between rows 45 and 46)
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D45") "0.317"
$ws.Range("E45").Value = "  -1.20%  "

# Row 46: coin identity/rank changed (values swapped between rows 45 and 46)
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D46") "3.27"
$ws.Range("E46").Value = "  -1.76%  "

$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("E48").Value = "  -4.57%  "
Set-TextValue $ws.Range("D49") "2.31"
$ws.Range("E49").Value = "  -3.91%  "
Set-TextValue $ws.Range("D50") "115.50"
$ws.Range("E50").Value = "  +6.26%  "
Set-TextValue $ws.Range("D51") "15.70"
$ws.Range("E51").Value = "  -4.36%  "
